$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.019.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.52%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.876.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'242.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.98%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4911"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.68%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2927"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.12%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06597"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.46%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.876.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.56%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'16.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.29%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07179"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.94%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.6653"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.32%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'86.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.895"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.85%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'29.983.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000007790"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -6.27%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.9984"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.03%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.57%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.122.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.30%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.9965"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.08%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.770"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.36%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.842"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.081"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.00%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.36%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'143.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +5.93%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'16.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.891"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -6.04%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.397"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.03%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.197"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.96%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08745"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.64%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.969"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.23%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05003"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.70%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7223"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.24%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.667"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.61%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.01820"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.40%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.95%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.149"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -7.07%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.9323"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.30%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.9990"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.753"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.4211"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.69%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'103.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.14%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'7.341"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -4.93%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.1269"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.11%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05697"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.17%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'32.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.267"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.3760"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.68%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.338"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.32%  "
$ws.Range("E51").Style = "Normal"
